$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# -----------------------------------------------------------------------
# 1) "base" command list (column F): insert the new command
#    assertMatch(text,regex) in alphabetical order (after assertEqual,
#    before assertNotContain), shifting the remaining entries down.
# -----------------------------------------------------------------------
$baseVals = @()
for ($r = 2; $r -le 42; $r++) {
    $baseVals += $ws.Cells.Item($r, 6).Value2
}
$newBaseVals = @($baseVals[0..8]) + @("assertMatch(text,regex)") + @($baseVals[9..($baseVals.Count - 1)])
for ($i = 0; $i -lt $newBaseVals.Count; $i++) {
    $ws.Cells.Item(2 + $i, 6).Value2 = $newBaseVals[$i]
}

# -----------------------------------------------------------------------
# 2) "external" command list (column J): insert the new command
#    openFile(filePath) in alphabetical order (before runJUnit), shifting
#    the remaining entries down.
# -----------------------------------------------------------------------
$externalVals = @()
for ($r = 2; $r -le 6; $r++) {
    $externalVals += $ws.Cells.Item($r, 10).Value2
}
$newExternalVals = @("openFile(filePath)") + $externalVals
for ($i = 0; $i -lt $newExternalVals.Count; $i++) {
    $ws.Cells.Item(2 + $i, 10).Value2 = $newExternalVals[$i]
}

# -----------------------------------------------------------------------
# 3) "target" category list (column A): the "tn.5250" category entry is
#    removed, shifting the following entries (web, webalert, webcookie,
#    ws, ws.async, xml) up by one row and clearing the now-unused last row.
# -----------------------------------------------------------------------
for ($r = 27; $r -le 32; $r++) {
    $v = $ws.Cells.Item($r + 1, 1).Value2
    $ws.Cells.Item($r, 1).Value2 = $v
}
$ws.Range("A33").ClearContents()

# -----------------------------------------------------------------------
# 4) The "tn.5250" data column (AA) is removed entirely; web/webalert/
#    webcookie/ws/ws.async/xml (previously AB:AG) each shift one column
#    to the left (AA:AF).
# -----------------------------------------------------------------------
$ws.Columns("AA").Delete()

# -----------------------------------------------------------------------
# 5) Update the workbook-level defined names to the new ranges.
# -----------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$F`$2:`$F`$45"
$wb.Names.Item("external").RefersTo = "='#system'!`$J`$2:`$J`$7"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$32"
$wb.Names.Item("web").RefersTo = "='#system'!`$AA`$2:`$AA`$151"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AC`$2:`$AC`$10"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AD`$2:`$AD`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AE`$2:`$AE`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AF`$2:`$AF`$27"

Write-Host "edit complete"
